$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "best tights",
    "best tights for women",
    "best workout pants",
    "best yoga leggings",
    "best yoga tops women",
    "better bodies leggings women",
    "better weigh",
    "bicycle distance tracker",
    "bicycle generator",
    "bicycle lift",
    "bicycle light generator",
    "bicycle reflective clothing",
    "bicycle spectrum",
    "bicycle support",
    "bicycle tights bike capris pants",
    "bicycle washing machine",
    "big and tall compression",
    "big and tall compression pants",
    "big and tall running pants",
    "big fig",
    "big game pro climbing sticks",
    "big joints",
    "big mens cold weather compression",
    "big muscle woman",
    "big trouble in thailand",
    "big win boxing",
    "bigger muscles",
    "bike brand coaches shorts",
    "bike capris",
    "bike coaches shorts",
    "bike coaches shorts polyester",
    "bike cold weather pants",
    "bike compression shorts women",
    "bike doctor",
    "bike generater",
    "bike generator",
    "bike hoop",
    "bike in place stand",
    "bike leg warmers",
    "bike mile tracker",
    "bike nook",
    "bike panties",
    "bike race pro",
    "bike race pro by top free games",
    "bike race pro free",
    "bike running shoes for women",
    "bike running shorts women",
    "bike shorts womens",
    "bike skirts for women",
    "bike swimmer",
    "bike tight",
    "bike tights",
    "bike tights women",
    "bike trick games",
    "bike wind pants",
    "bike winter gear women",
    "biking capris",
    "biking gear",
    "biking knee brace",
    "biking knee warmers",
    "biking oants",
    "biking tights women",
    "bills gear for women",
    "bills leggings",
    "bioflect compression leggings",
    "bioflex active",
    "bit trip runner 3",
    "black 3",
    "black and green compression pants men",
    "black and mild cigarettes",
    "black and pink tights",
    "black and red stripe leggings for women",
    "black and red tights",
    "black and red tights for women",
    "black and white boxing shorts",
    "black and white high waisted leggings",
    "black and white tights for women",
    "black apparel",
    "black apparel fabric",
    "black athletic",
    "black athletic capris for women",
    "black athletic leggings",
    "black athletic leggings for women",
    "black athletic pants",
    "black athletic pants women",
    "black athletic tape",
    "black ballet leggings",
    "black ballet tights",
    "black basketball leggings for men",
    "black boys compression pants",
    "black capri legging",
    "black capri leggings",
    "black capri leggings for women",
    "black capri pants",
    "black capri women",
    "black capri workout leggings",
    "black capri yoga pants",
    "black capris",
    "black capris for women",
    "black capris leggings"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
